# Apply the round-robin schedule fix: rotate home/away assignments
# (circle method, round 1 untouched order) and align each player's
# rating value to travel with their name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 2).Value = "Петров"
$ws.Cells.Item(3, 3).Value = 840
$ws.Cells.Item(3, 5).Value = "Котов"
$ws.Cells.Item(3, 6).Value = 900
# Row 4
$ws.Cells.Item(4, 2).Value = "Серов"
$ws.Cells.Item(4, 3).Value = 730
$ws.Cells.Item(4, 5).Value = "Амелин"
$ws.Cells.Item(4, 6).Value = 680
# Row 5
$ws.Cells.Item(5, 2).Value = "Белов"
$ws.Cells.Item(5, 3).Value = 810
$ws.Cells.Item(5, 5).Value = "Сидоров"
$ws.Cells.Item(5, 6).Value = 720
# Row 6
$ws.Cells.Item(6, 2).Value = "Кротов"
$ws.Cells.Item(6, 3).Value = 800
$ws.Cells.Item(6, 5).Value = "пропуск"
$ws.Cells.Item(6, 6).Value = 0
# Row 7
$ws.Cells.Item(7, 5).Value = "Перов"
$ws.Cells.Item(7, 6).Value = 750
# Row 8
$ws.Cells.Item(8, 2).Value = "Якин"
$ws.Cells.Item(8, 3).Value = 750
$ws.Cells.Item(8, 5).Value = "Уткин"
$ws.Cells.Item(8, 6).Value = 750
# Row 10
$ws.Cells.Item(10, 2).Value = "Котов"
$ws.Cells.Item(10, 3).Value = 900
$ws.Cells.Item(10, 5).Value = "Уткин"
$ws.Cells.Item(10, 6).Value = 750
# Row 11
$ws.Cells.Item(11, 2).Value = "Перов"
$ws.Cells.Item(11, 3).Value = 750
$ws.Cells.Item(11, 5).Value = "Якин"
$ws.Cells.Item(11, 6).Value = 750
# Row 12
$ws.Cells.Item(12, 2).Value = "пропуск"
$ws.Cells.Item(12, 3).Value = 0
# Row 13
$ws.Cells.Item(13, 2).Value = "Сидоров"
$ws.Cells.Item(13, 3).Value = 720
$ws.Cells.Item(13, 5).Value = "Кротов"
$ws.Cells.Item(13, 6).Value = 800
# Row 14
$ws.Cells.Item(14, 2).Value = "Амелин"
$ws.Cells.Item(14, 3).Value = 680
$ws.Cells.Item(14, 5).Value = "Белов"
$ws.Cells.Item(14, 6).Value = 810
# Row 15
$ws.Cells.Item(15, 2).Value = "Петров"
$ws.Cells.Item(15, 3).Value = 840
$ws.Cells.Item(15, 5).Value = "Серов"
$ws.Cells.Item(15, 6).Value = 730
# Row 17
$ws.Cells.Item(17, 2).Value = "Серов"
$ws.Cells.Item(17, 3).Value = 730
$ws.Cells.Item(17, 5).Value = "Котов"
$ws.Cells.Item(17, 6).Value = 900
# Row 18
$ws.Cells.Item(18, 2).Value = "Белов"
$ws.Cells.Item(18, 3).Value = 810
$ws.Cells.Item(18, 5).Value = "Петров"
$ws.Cells.Item(18, 6).Value = 840
# Row 19
$ws.Cells.Item(19, 2).Value = "Кротов"
$ws.Cells.Item(19, 3).Value = 800
$ws.Cells.Item(19, 5).Value = "Амелин"
$ws.Cells.Item(19, 6).Value = 680
# Row 20
$ws.Cells.Item(20, 5).Value = "Сидоров"
$ws.Cells.Item(20, 6).Value = 720
# Row 21
$ws.Cells.Item(21, 2).Value = "Якин"
$ws.Cells.Item(21, 3).Value = 750
$ws.Cells.Item(21, 5).Value = "пропуск"
$ws.Cells.Item(21, 6).Value = 0
# Row 22
$ws.Cells.Item(22, 2).Value = "Уткин"
$ws.Cells.Item(22, 3).Value = 750
$ws.Cells.Item(22, 5).Value = "Перов"
$ws.Cells.Item(22, 6).Value = 750
# Row 24
$ws.Cells.Item(24, 2).Value = "Котов"
$ws.Cells.Item(24, 3).Value = 900
$ws.Cells.Item(24, 5).Value = "Перов"
$ws.Cells.Item(24, 6).Value = 750
# Row 25
$ws.Cells.Item(25, 2).Value = "пропуск"
$ws.Cells.Item(25, 3).Value = 0
$ws.Cells.Item(25, 5).Value = "Уткин"
$ws.Cells.Item(25, 6).Value = 750
# Row 26
$ws.Cells.Item(26, 2).Value = "Сидоров"
$ws.Cells.Item(26, 3).Value = 720
$ws.Cells.Item(26, 5).Value = "Якин"
$ws.Cells.Item(26, 6).Value = 750
# Row 27
$ws.Cells.Item(27, 2).Value = "Амелин"
$ws.Cells.Item(27, 3).Value = 680
# Row 28
$ws.Cells.Item(28, 2).Value = "Петров"
$ws.Cells.Item(28, 3).Value = 840
$ws.Cells.Item(28, 5).Value = "Кротов"
$ws.Cells.Item(28, 6).Value = 800
# Row 29
$ws.Cells.Item(29, 2).Value = "Серов"
$ws.Cells.Item(29, 3).Value = 730
$ws.Cells.Item(29, 5).Value = "Белов"
$ws.Cells.Item(29, 6).Value = 810
# Row 31
$ws.Cells.Item(31, 2).Value = "Белов"
$ws.Cells.Item(31, 3).Value = 810
$ws.Cells.Item(31, 5).Value = "Котов"
$ws.Cells.Item(31, 6).Value = 900
# Row 32
$ws.Cells.Item(32, 2).Value = "Кротов"
$ws.Cells.Item(32, 3).Value = 800
$ws.Cells.Item(32, 5).Value = "Серов"
$ws.Cells.Item(32, 6).Value = 730
# Row 33
$ws.Cells.Item(33, 5).Value = "Петров"
$ws.Cells.Item(33, 6).Value = 840
# Row 34
$ws.Cells.Item(34, 2).Value = "Якин"
$ws.Cells.Item(34, 3).Value = 750
$ws.Cells.Item(34, 5).Value = "Амелин"
$ws.Cells.Item(34, 6).Value = 680
# Row 35
$ws.Cells.Item(35, 2).Value = "Уткин"
$ws.Cells.Item(35, 3).Value = 750
$ws.Cells.Item(35, 5).Value = "Сидоров"
$ws.Cells.Item(35, 6).Value = 720
# Row 36
$ws.Cells.Item(36, 2).Value = "Перов"
$ws.Cells.Item(36, 3).Value = 750
$ws.Cells.Item(36, 5).Value = "пропуск"
$ws.Cells.Item(36, 6).Value = 0
# Row 38
$ws.Cells.Item(38, 2).Value = "Котов"
$ws.Cells.Item(38, 3).Value = 900
$ws.Cells.Item(38, 5).Value = "пропуск"
$ws.Cells.Item(38, 6).Value = 0
# Row 39
$ws.Cells.Item(39, 2).Value = "Сидоров"
$ws.Cells.Item(39, 3).Value = 720
$ws.Cells.Item(39, 5).Value = "Перов"
$ws.Cells.Item(39, 6).Value = 750
# Row 40
$ws.Cells.Item(40, 2).Value = "Амелин"
$ws.Cells.Item(40, 3).Value = 680
$ws.Cells.Item(40, 5).Value = "Уткин"
$ws.Cells.Item(40, 6).Value = 750
# Row 41
$ws.Cells.Item(41, 2).Value = "Петров"
$ws.Cells.Item(41, 3).Value = 840
$ws.Cells.Item(41, 5).Value = "Якин"
$ws.Cells.Item(41, 6).Value = 750
# Row 42
$ws.Cells.Item(42, 2).Value = "Серов"
$ws.Cells.Item(42, 3).Value = 730
# Row 43
$ws.Cells.Item(43, 2).Value = "Белов"
$ws.Cells.Item(43, 3).Value = 810
$ws.Cells.Item(43, 5).Value = "Кротов"
$ws.Cells.Item(43, 6).Value = 800
# Row 45
$ws.Cells.Item(45, 2).Value = "Кротов"
$ws.Cells.Item(45, 3).Value = 800
$ws.Cells.Item(45, 5).Value = "Котов"
$ws.Cells.Item(45, 6).Value = 900
# Row 46
$ws.Cells.Item(46, 5).Value = "Белов"
$ws.Cells.Item(46, 6).Value = 810
# Row 47
$ws.Cells.Item(47, 2).Value = "Якин"
$ws.Cells.Item(47, 3).Value = 750
$ws.Cells.Item(47, 5).Value = "Серов"
$ws.Cells.Item(47, 6).Value = 730
# Row 48
$ws.Cells.Item(48, 2).Value = "Уткин"
$ws.Cells.Item(48, 3).Value = 750
$ws.Cells.Item(48, 5).Value = "Петров"
$ws.Cells.Item(48, 6).Value = 840
# Row 49
$ws.Cells.Item(49, 2).Value = "Перов"
$ws.Cells.Item(49, 3).Value = 750
$ws.Cells.Item(49, 5).Value = "Амелин"
$ws.Cells.Item(49, 6).Value = 680
# Row 50
$ws.Cells.Item(50, 2).Value = "пропуск"
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 5).Value = "Сидоров"
$ws.Cells.Item(50, 6).Value = 720
# Row 52
$ws.Cells.Item(52, 2).Value = "Котов"
$ws.Cells.Item(52, 3).Value = 900
$ws.Cells.Item(52, 5).Value = "Сидоров"
$ws.Cells.Item(52, 6).Value = 720
# Row 53
$ws.Cells.Item(53, 2).Value = "Амелин"
$ws.Cells.Item(53, 3).Value = 680
$ws.Cells.Item(53, 5).Value = "пропуск"
$ws.Cells.Item(53, 6).Value = 0
# Row 54
$ws.Cells.Item(54, 2).Value = "Петров"
$ws.Cells.Item(54, 3).Value = 840
$ws.Cells.Item(54, 5).Value = "Перов"
$ws.Cells.Item(54, 6).Value = 750
# Row 55
$ws.Cells.Item(55, 2).Value = "Серов"
$ws.Cells.Item(55, 3).Value = 730
$ws.Cells.Item(55, 5).Value = "Уткин"
$ws.Cells.Item(55, 6).Value = 750
# Row 56
$ws.Cells.Item(56, 2).Value = "Белов"
$ws.Cells.Item(56, 3).Value = 810
$ws.Cells.Item(56, 5).Value = "Якин"
$ws.Cells.Item(56, 6).Value = 750
# Row 57
$ws.Cells.Item(57, 2).Value = "Кротов"
$ws.Cells.Item(57, 3).Value = 800
# Row 59
$ws.Cells.Item(59, 5).Value = "Котов"
$ws.Cells.Item(59, 6).Value = 900
# Row 60
$ws.Cells.Item(60, 2).Value = "Якин"
$ws.Cells.Item(60, 3).Value = 750
$ws.Cells.Item(60, 5).Value = "Кротов"
$ws.Cells.Item(60, 6).Value = 800
# Row 61
$ws.Cells.Item(61, 2).Value = "Уткин"
$ws.Cells.Item(61, 3).Value = 750
$ws.Cells.Item(61, 5).Value = "Белов"
$ws.Cells.Item(61, 6).Value = 810
# Row 62
$ws.Cells.Item(62, 2).Value = "Перов"
$ws.Cells.Item(62, 3).Value = 750
$ws.Cells.Item(62, 5).Value = "Серов"
$ws.Cells.Item(62, 6).Value = 730
# Row 63
$ws.Cells.Item(63, 2).Value = "пропуск"
$ws.Cells.Item(63, 3).Value = 0
$ws.Cells.Item(63, 5).Value = "Петров"
$ws.Cells.Item(63, 6).Value = 840
# Row 64
$ws.Cells.Item(64, 2).Value = "Сидоров"
$ws.Cells.Item(64, 3).Value = 720
$ws.Cells.Item(64, 5).Value = "Амелин"
$ws.Cells.Item(64, 6).Value = 680
# Row 66
$ws.Cells.Item(66, 2).Value = "Котов"
$ws.Cells.Item(66, 3).Value = 900
$ws.Cells.Item(66, 5).Value = "Амелин"
$ws.Cells.Item(66, 6).Value = 680
# Row 67
$ws.Cells.Item(67, 2).Value = "Петров"
$ws.Cells.Item(67, 3).Value = 840
$ws.Cells.Item(67, 5).Value = "Сидоров"
$ws.Cells.Item(67, 6).Value = 720
# Row 68
$ws.Cells.Item(68, 2).Value = "Серов"
$ws.Cells.Item(68, 3).Value = 730
$ws.Cells.Item(68, 5).Value = "пропуск"
$ws.Cells.Item(68, 6).Value = 0
# Row 69
$ws.Cells.Item(69, 2).Value = "Белов"
$ws.Cells.Item(69, 3).Value = 810
$ws.Cells.Item(69, 5).Value = "Перов"
$ws.Cells.Item(69, 6).Value = 750
# Row 70
$ws.Cells.Item(70, 2).Value = "Кротов"
$ws.Cells.Item(70, 3).Value = 800
$ws.Cells.Item(70, 5).Value = "Уткин"
$ws.Cells.Item(70, 6).Value = 750
# Row 71
$ws.Cells.Item(71, 5).Value = "Якин"
$ws.Cells.Item(71, 6).Value = 750
# Row 73
$ws.Cells.Item(73, 2).Value = "Якин"
$ws.Cells.Item(73, 3).Value = 750
$ws.Cells.Item(73, 5).Value = "Котов"
$ws.Cells.Item(73, 6).Value = 900
# Row 74
$ws.Cells.Item(74, 2).Value = "Уткин"
$ws.Cells.Item(74, 3).Value = 750
# Row 75
$ws.Cells.Item(75, 2).Value = "Перов"
$ws.Cells.Item(75, 3).Value = 750
$ws.Cells.Item(75, 5).Value = "Кротов"
$ws.Cells.Item(75, 6).Value = 800
# Row 76
$ws.Cells.Item(76, 2).Value = "пропуск"
$ws.Cells.Item(76, 3).Value = 0
$ws.Cells.Item(76, 5).Value = "Белов"
$ws.Cells.Item(76, 6).Value = 810
# Row 77
$ws.Cells.Item(77, 2).Value = "Сидоров"
$ws.Cells.Item(77, 3).Value = 720
$ws.Cells.Item(77, 5).Value = "Серов"
$ws.Cells.Item(77, 6).Value = 730
# Row 78
$ws.Cells.Item(78, 2).Value = "Амелин"
$ws.Cells.Item(78, 3).Value = 680
$ws.Cells.Item(78, 5).Value = "Петров"
$ws.Cells.Item(78, 6).Value = 840
